# Insert a new weekly price record into the Orégano price history sheet.
# A new row is inserted at row 222, pushing the existing rows 222-298 down
# to 223-299 (the rest of the table is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 222.
$ws.Rows.Item(222).Insert()

# Populate the new row with the latest weekly record.
$ws.Cells.Item(222, 1).Value = 6
$ws.Cells.Item(222, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(222, 3).Value = "Metropolitana"
$ws.Cells.Item(222, 4).Value = 44985
$ws.Cells.Item(222, 5).Value = 13
$ws.Cells.Item(222, 6).Value = 100112029
$ws.Cells.Item(222, 7).Value = "Orégano"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 48
$ws.Cells.Item(222, 11).Value = 16000
$ws.Cells.Item(222, 12).Value = 17000
$ws.Cells.Item(222, 13).Value = 16458
$ws.Cells.Item(222, 14).Value = "`$/docena de atados"
$ws.Cells.Item(222, 15).Value = "Región Metropolitana"
$ws.Cells.Item(222, 16).Value = 5486
$ws.Cells.Item(222, 17).Value = 3
$ws.Cells.Item(222, 18).Value = "Hortaliza"
